# Session2.pptx - title slide update
#
# The only substantive content edit in the target revision is on the title
# slide (slide 1), in the "CustomShape 2" textbox that holds the presenter /
# term byline under the title. Its second line changes from "Hilary 2022" to
# "Trinity 2022" - i.e. the term name "Hilary" is swapped for "Trinity" while
# the year stays the same, and PowerPoint ends up recording that line as two
# separate runs ("Trinity " and "2022") instead of one.
#
# We reproduce that by grabbing just the "Hilary " substring (the word plus
# the following space) via TextRange.Characters and overwriting its .Text -
# this is exactly what PowerPoint does when a user selects part of a run and
# retypes it: the paragraph is split into a run for the edited span and a run
# for the untouched remainder, both keeping the original character formatting
# (Calibri 32pt grey 8B8B8B).

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(1)
$sh = $s.Shapes.Item(2)          # "CustomShape 2" - byline textbox
$tr = $sh.TextFrame.TextRange

$byline = $tr.Paragraphs(2)      # paragraph 2 = "Hilary 2022"

$termWord = $byline.Characters(1, 7)   # "Hilary " (word + trailing space)
$termWord.Text = "Trinity "

Write-Output $tr.Paragraphs(1).Text
Write-Output $tr.Paragraphs(2).Text
